{"js": "// Update the five filled rows of the division-problems table.\n// Row 3 (the \"30\u00f79=\" row) also loses its first cell and gains a new\n// cell at the end; since every cell in the table shares identical\n// formatting (tcPr/pPr/rPr), rewriting that row's five cell values in\n// place (\"53\u00f75=,13\u00f79=,14\u00f76=,32\u00f75=,52\u00f76=\") reproduces the same result as\n// deleting the old first cell and appending a freshly formatted one.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices (within the 20-row table) that contain the visible problems.\nconst rowIndexes = [0, 4, 8, 12, 16];\n\nconst newRowValues = [\n  [\"92\u00f72=\", \"25\u00f74=\", \"36\u00f73=\", \"93\u00f75=\", \"47\u00f76=\"],\n  [\"52\u00f73=\", \"65\u00f76=\", \"97\u00f74=\", \"40\u00f76=\", \"45\u00f72=\"],\n  [\"53\u00f75=\", \"13\u00f79=\", \"14\u00f76=\", \"32\u00f75=\", \"52\u00f76=\"],\n  [\"50\u00f72=\", \"78\u00f73=\", \"75\u00f74=\", \"45\u00f75=\", \"85\u00f73=\"],\n  [\"18\u00f78=\", \"22\u00f79=\", \"24\u00f76=\", \"10\u00f76=\", \"62\u00f78=\"],\n];\n\nfor (let i = 0; i < rowIndexes.length; i++) {\n  const r = rowIndexes[i];\n  const values = newRowValues[i];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the five filled rows of the division-problems table.\n# Row 9 (the \"30\u00f79=\" row) also loses its first cell and gains a new\n# cell at the end; since every cell in the table shares identical\n# formatting, rewriting that row's five cell values in place\n# (\"53\u00f75=,13\u00f79=,14\u00f76=,32\u00f75=,52\u00f76=\") reproduces the same result as\n# deleting the old first cell and appending a freshly formatted one.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowNumbers = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n    @(\"92\u00f72=\", \"25\u00f74=\", \"36\u00f73=\", \"93\u00f75=\", \"47\u00f76=\"),\n    @(\"52\u00f73=\", \"65\u00f76=\", \"97\u00f74=\", \"40\u00f76=\", \"45\u00f72=\"),\n    @(\"53\u00f75=\", \"13\u00f79=\", \"14\u00f76=\", \"32\u00f75=\", \"52\u00f76=\"),\n    @(\"50\u00f72=\", \"78\u00f73=\", \"75\u00f74=\", \"45\u00f75=\", \"85\u00f73=\"),\n    @(\"18\u00f78=\", \"22\u00f79=\", \"24\u00f76=\", \"10\u00f76=\", \"62\u00f78=\")\n)\n\nfor ($i = 0; $i -lt $rowNumbers.Length; $i++) {\n    $r = $rowNumbers[$i]\n    $values = $newValues[$i]\n    for ($c = 1; $c -le $values.Length; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$c - 1]\n    }\n}\n"}
